$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the value in A2 (change link/name in index)
$ws.Range("A2").Value = "iiandjdmd"

# Remove rows 3 and 4 entirely (they duplicated the old data)
$ws.Rows("3:4").Delete()
